# Update the "Pais" (countries) sheet with refreshed COVID-19 stats and a
# newer "last updated" timestamp, and fix the Togo/Malta row mix-up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Refresh the "last updated" timestamp (row 1) -----------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Agosto de 2020 a las 13:28"

# --- 2. Fix Togo / Malta: row 151 used to show Togo's figures and row 152
#        Malta's; the countries swap rows and Malta's figures are refreshed.
$ws.Cells.Item(151, 1).Value = "Malta"
$ws.Cells.Item(151, 2).Value = 1089
$ws.Cells.Item(151, 3).Value = 54
$ws.Cells.Item(151, 4).Value = 684
$ws.Cells.Item(151, 5).Value = 396
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 9

$ws.Cells.Item(152, 1).Value = "Togo"
$ws.Cells.Item(152, 2).Value = 1046
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(152, 4).Value = 721
$ws.Cells.Item(152, 5).Value = 302
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 23

# --- 3. Refresh per-country case counts -------------------------------------
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes
$rows        = @(14,     36,    37,    42,    43,    79,    158)   # sheet row
$countries   = @("Iran", "Oman","Ucrania","Bielorrusia","Emiratos Arabes Unidos","Estado de Palestina","Vietnam")
$totales     = @(326712, 81580, 80949, 68850, 62525, 14205, 841)
$nuevos      = @(2020,   223,   1199,  112,   225,   277,   29)
$activos     = @(284371, 74691, 43972, 64935, 56568, 7945,  395)
$recuperados = @(23914,  6376,  35080, 3328,  5600,  6164,  435)
$criticos    = @(0,      0,     0,     0,     0,     0,     0)
$muertesHoy  = @(163,    4,     18,    2,     1,     0,     1)
$muertes     = @(18427,  513,   1897,  587,   357,   96,    11)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $totales[$i]
    $ws.Cells.Item($r, 3).Value = $nuevos[$i]
    $ws.Cells.Item($r, 4).Value = $activos[$i]
    $ws.Cells.Item($r, 5).Value = $recuperados[$i]
    $ws.Cells.Item($r, 6).Value = $criticos[$i]
    $ws.Cells.Item($r, 7).Value = $muertesHoy[$i]
    $ws.Cells.Item($r, 8).Value = $muertes[$i]
}
